# Update the equipment name labels that are shared between both sheets
# (sharedStrings entries are reused across "复杂度1" and "复杂度2" sheets, as
# well as in the embedded charts' category caches).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("复杂度1")
$ws2 = $wb.Worksheets.Item("复杂度2")

$ws1.Range("A2").Value = "消防无人机"
$ws1.Range("A3").Value = "植保无人机"
$ws1.Range("A4").Value = "电力巡检机器人"

$ws2.Range("A2").Value = "消防无人机"
$ws2.Range("A3").Value = "植保无人机"
$ws2.Range("A4").Value = "电力巡检机器人"

# Restore the originally selected cells on each sheet
$ws1.Range("A4").Select()
$ws2.Range("B26").Select()

# Make "复杂度2" the active sheet/tab (it was already the active tab)
$ws2.Activate()

# Widen the workbook window
$excel.ActiveWindow.Width = 24750
